# Update 556 barrel extensions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("556-muzzles")

# Row 26 - IPS 4 (ips_5.56x45_4inch_barrel_extension)
$ws.Range("M26").Value = -3
$ws.Range("O26").Value = -1
$ws.Range("P26").Value = -1
$ws.Range("Q26").Value = 0.2
$ws.Range("R26").Value = 0.04

# Row 27 - IPS 2 (ips_5.56x45_2inch_barrel_extension)
$ws.Range("M27").Value = -2
$ws.Range("O27").Value = 0
$ws.Range("Q27").Value = 0.1
$ws.Range("R27").Value = 0.02

# Row 28 - IPS 1 (ips_5.56x45_1inch_barrel_extension)
$ws.Range("M28").Value = -1
$ws.Range("Q28").Value = 0.05
$ws.Range("R28").Value = 0.01

# Update view state to match final selection/scroll position
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M29").Select()
